$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
